$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 10.105724
$ws.Cells.Item(2, 8).Value = 30.317172
$ws.Cells.Item(2, 9).Value = 0.5504853801993582
$ws.Cells.Item(2, 10).Value = 0.5504853801993582
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 409.6166503333334
$ws.Cells.Item(2, 14).Value = 1228.849951
$ws.Cells.Item(2, 15).Value = 0.6234125531262766
$ws.Cells.Item(2, 16).Value = 0.6234125531262766
$ws.Cells.Item(2, 17).Value = 4139.472814073175
$ws.Cells.Item(2, 18).Value = 37255.25532665857
$ws.Cells.Item(2, 19).Value = 0.343179496328771
$ws.Cells.Item(2, 20).Value = 0.343179496328771

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 10.105724
$ws.Cells.Item(3, 8).Value = 30.317172
$ws.Cells.Item(3, 9).Value = 0.5504853801993582
$ws.Cells.Item(3, 10).Value = 0.5504853801993582
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 56.495384
$ws.Cells.Item(3, 14).Value = 169.486152
$ws.Cells.Item(3, 15).Value = 0.08598266586728959
$ws.Cells.Item(3, 16).Value = 0.08598266586728959
$ws.Cells.Item(3, 17).Value = 570.926757978016
$ws.Cells.Item(3, 18).Value = 5138.340821802144
$ws.Cells.Item(3, 19).Value = 0.04733220051050929
$ws.Cells.Item(3, 20).Value = 0.04733220051050929

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 10.105724
$ws.Cells.Item(4, 8).Value = 30.317172
$ws.Cells.Item(4, 9).Value = 0.5504853801993582
$ws.Cells.Item(4, 10).Value = 0.5504853801993582
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 190.9434713333333
$ws.Cells.Item(4, 14).Value = 572.830414
$ws.Cells.Item(4, 15).Value = 0.2906047810064339
$ws.Cells.Item(4, 16).Value = 0.2906047810064338
$ws.Cells.Item(4, 17).Value = 1929.622020896579
$ws.Cells.Item(4, 18).Value = 17366.59818806921
$ws.Cells.Item(4, 19).Value = 0.159973683360078
$ws.Cells.Item(4, 20).Value = 0.159973683360078

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 5.009378000000001
$ws.Cells.Item(5, 8).Value = 15.028134
$ws.Cells.Item(5, 9).Value = 0.2728740021884925
$ws.Cells.Item(5, 10).Value = 0.2728740021884924
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 409.6166503333334
$ws.Cells.Item(5, 14).Value = 1228.849951
$ws.Cells.Item(5, 15).Value = 0.6234125531262766
$ws.Cells.Item(5, 16).Value = 0.6234125531262766
$ws.Cells.Item(5, 17).Value = 2051.924636613493
$ws.Cells.Item(5, 18).Value = 18467.32172952144
$ws.Cells.Item(5, 19).Value = 0.1701130783861133
$ws.Cells.Item(5, 20).Value = 0.1701130783861133

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5.009378000000001
$ws.Cells.Item(6, 8).Value = 15.028134
$ws.Cells.Item(6, 9).Value = 0.2728740021884925
$ws.Cells.Item(6, 10).Value = 0.2728740021884924
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 56.495384
$ws.Cells.Item(6, 14).Value = 169.486152
$ws.Cells.Item(6, 15).Value = 0.08598266586728959
$ws.Cells.Item(6, 16).Value = 0.08598266586728959
$ws.Cells.Item(6, 17).Value = 283.006733711152
$ws.Cells.Item(6, 18).Value = 2547.060603400368
$ws.Cells.Item(6, 19).Value = 0.0234624341540432
$ws.Cells.Item(6, 20).Value = 0.02346243415404319

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5.009378000000001
$ws.Cells.Item(7, 8).Value = 15.028134
$ws.Cells.Item(7, 9).Value = 0.2728740021884925
$ws.Cells.Item(7, 10).Value = 0.2728740021884924
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 190.9434713333333
$ws.Cells.Item(7, 14).Value = 572.830414
$ws.Cells.Item(7, 15).Value = 0.2906047810064339
$ws.Cells.Item(7, 16).Value = 0.2906047810064338
$ws.Cells.Item(7, 17).Value = 956.5080245408309
$ws.Cells.Item(7, 18).Value = 8608.572220867478
$ws.Cells.Item(7, 19).Value = 0.07929848964833602
$ws.Cells.Item(7, 20).Value = 0.079298489648336

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 3.242740666666667
$ws.Cells.Item(8, 8).Value = 9.728222
$ws.Cells.Item(8, 9).Value = 0.1766406176121494
$ws.Cells.Item(8, 10).Value = 0.1766406176121493
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 409.6166503333334
$ws.Cells.Item(8, 14).Value = 1228.849951
$ws.Cells.Item(8, 15).Value = 0.6234125531262766
$ws.Cells.Item(8, 16).Value = 0.6234125531262766
$ws.Cells.Item(8, 17).Value = 1328.28056977968
$ws.Cells.Item(8, 18).Value = 11954.52512801712
$ws.Cells.Item(8, 19).Value = 0.1101199784113924
$ws.Cells.Item(8, 20).Value = 0.1101199784113924

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 3.242740666666667
$ws.Cells.Item(9, 8).Value = 9.728222
$ws.Cells.Item(9, 9).Value = 0.1766406176121494
$ws.Cells.Item(9, 10).Value = 0.1766406176121493
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 56.495384
$ws.Cells.Item(9, 14).Value = 169.486152
$ws.Cells.Item(9, 15).Value = 0.08598266586728959
$ws.Cells.Item(9, 16).Value = 0.08598266586728959
$ws.Cells.Item(9, 17).Value = 183.1998791757493
$ws.Cells.Item(9, 18).Value = 1648.798912581744
$ws.Cells.Item(9, 19).Value = 0.01518803120273711
$ws.Cells.Item(9, 20).Value = 0.0151880312027371

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 3.242740666666667
$ws.Cells.Item(10, 8).Value = 9.728222
$ws.Cells.Item(10, 9).Value = 0.1766406176121494
$ws.Cells.Item(10, 10).Value = 0.1766406176121493
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 190.9434713333333
$ws.Cells.Item(10, 14).Value = 572.830414
$ws.Cells.Item(10, 15).Value = 0.2906047810064339
$ws.Cells.Item(10, 16).Value = 0.2906047810064338
$ws.Cells.Item(10, 17).Value = 619.180159527101
$ws.Cells.Item(10, 18).Value = 5572.621435743908
$ws.Cells.Item(10, 19).Value = 0.05133260799801989
$ws.Cells.Item(10, 20).Value = 0.05133260799801988
